# Natmi following Dr Hou advice
# Update ligand/receptor-expressing cell counts (E, K) from 1 to 3,
# and recompute the dependent expression/specificity metrics for rows 2-19.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    2 = @{ "E" = 3; "G" = 2.663313; "H" = 7.989939000000001; "I" = 0.3794306644527501; "J" = 0.3794306644527502; "K" = 3; "M" = 34.218763; "N" = 102.656289; "O" = 0.4046921425624349; "P" = 0.4046921425624349; "Q" = 91.13527634181901; "R" = 820.2174870763711; "S" = 0.1535526085512718; "T" = 0.1535526085512718 }
    3 = @{ "E" = 3; "G" = 2.663313; "H" = 7.989939000000001; "I" = 0.3794306644527501; "J" = 0.3794306644527502; "K" = 3; "M" = 13.95683033333333; "N" = 41.870491; "O" = 0.165062061740135; "P" = 0.165062061740135; "Q" = 37.17140766556101; "R" = 334.542668990049; "S" = 0.06262960776200031; "T" = 0.06262960776200031 }
    4 = @{ "E" = 3; "G" = 2.663313; "H" = 7.989939000000001; "I" = 0.3794306644527501; "J" = 0.3794306644527502; "K" = 3; "M" = 14.45863566666667; "N" = 43.375907; "O" = 0.1709967203219173; "P" = 0.1709967203219173; "Q" = 38.507872333297; "R" = 346.570850999673; "S" = 0.06488139921098617; "T" = 0.06488139921098617 }
    5 = @{ "E" = 3; "G" = 2.663313; "H" = 7.989939000000001; "I" = 0.3794306644527501; "J" = 0.3794306644527502; "K" = 3; "M" = 10.640006; "N" = 31.920018; "O" = 0.1258352566694817; "P" = 0.1258352566694817; "Q" = 28.33766629987801; "R" = 255.038996698902; "S" = 0.04774575504968379; "T" = 0.0477457550496838 }
    6 = @{ "E" = 3; "G" = 2.663313; "H" = 7.989939000000001; "I" = 0.3794306644527501; "J" = 0.3794306644527502; "K" = 3; "M" = 8.185362; "N" = 24.556086; "O" = 0.09680512663269379; "P" = 0.09680512663269379; "Q" = 21.800181024306; "R" = 196.201629218754; "S" = 0.03673083352067563; "T" = 0.03673083352067563 }
    7 = @{ "E" = 3; "G" = 2.663313; "H" = 7.989939000000001; "I" = 0.3794306644527501; "J" = 0.3794306644527502; "K" = 3; "M" = 3.095449666666667; "N" = 9.286349; "O" = 0.03660869207333731; "P" = 0.03660869207333731; "Q" = 8.244151338079; "R" = 74.19736204271101; "S" = 0.0138904603581325; "T" = 0.01389046035813251 }
    8 = @{ "E" = 3; "G" = 3.178631333333334; "H" = 9.535894000000001; "I" = 0.4528458348143826; "J" = 0.4528458348143827; "K" = 3; "M" = 34.218763; "N" = 102.656289; "O" = 0.4046921425624349; "P" = 0.4046921425624349; "Q" = 108.7688322597074; "R" = 978.9194903373661; "S" = 0.183263151141507; "T" = 0.183263151141507 }
    9 = @{ "E" = 3; "G" = 3.178631333333334; "H" = 9.535894000000001; "I" = 0.4528458348143826; "J" = 0.4528458348143827; "K" = 3; "M" = 13.95683033333333; "N" = 41.870491; "O" = 0.165062061740135; "P" = 0.165062061740135; "Q" = 44.36361821155045; "R" = 399.2725639039541; "S" = 0.07474766714489461; "T" = 0.07474766714489463 }
    10 = @{ "E" = 3; "G" = 3.178631333333334; "H" = 9.535894000000001; "I" = 0.4528458348143826; "J" = 0.4528458348143827; "K" = 3; "M" = 14.45863566666667; "N" = 43.375907; "O" = 0.1709967203219173; "P" = 0.1709967203219173; "Q" = 45.95867236731755; "R" = 413.628051305858; "S" = 0.07743515256470014; "T" = 0.07743515256470015 }
    11 = @{ "E" = 3; "G" = 3.178631333333334; "H" = 9.535894000000001; "I" = 0.4528458348143826; "J" = 0.4528458348143827; "K" = 3; "M" = 10.640006; "N" = 31.920018; "O" = 0.1258352566694817; "P" = 0.1258352566694817; "Q" = 33.82065645845467; "R" = 304.385908126092; "S" = 0.05698397185557354; "T" = 0.05698397185557354 }
    12 = @{ "E" = 3; "G" = 3.178631333333334; "H" = 9.535894000000001; "I" = 0.4528458348143826; "J" = 0.4528458348143827; "K" = 3; "M" = 8.185362; "N" = 24.556086; "O" = 0.09680512663269379; "P" = 0.09680512663269379; "Q" = 26.018248127876; "R" = 234.164233150884; "S" = 0.04383779838429425; "T" = 0.04383779838429425 }
    13 = @{ "E" = 3; "G" = 3.178631333333334; "H" = 9.535894000000001; "I" = 0.4528458348143826; "J" = 0.4528458348143827; "K" = 3; "M" = 3.095449666666667; "N" = 9.286349; "O" = 0.03660869207333731; "P" = 0.03660869207333731; "Q" = 9.83929330122289; "R" = 88.553639711006; "S" = 0.01657809372341311; "T" = 0.01657809372341311 }
    14 = @{ "E" = 3; "G" = 1.177290666666667; "H" = 3.531872; "I" = 0.1677235007328671; "J" = 0.1677235007328671; "K" = 3; "M" = 34.218763; "N" = 102.656289; "O" = 0.4046921425624349; "P" = 0.4046921425624349; "Q" = 40.28543030477867; "R" = 362.5688727430081; "S" = 0.06787638286965612; "T" = 0.06787638286965611 }
    15 = @{ "E" = 3; "G" = 1.177290666666667; "H" = 3.531872; "I" = 0.1677235007328671; "J" = 0.1677235007328671; "K" = 3; "M" = 13.95683033333333; "N" = 41.870491; "O" = 0.165062061740135; "P" = 0.165062061740135; "Q" = 16.43124608768356; "R" = 147.881214789152; "S" = 0.0276847868332401; "T" = 0.0276847868332401 }
    16 = @{ "E" = 3; "G" = 1.177290666666667; "H" = 3.531872; "I" = 0.1677235007328671; "J" = 0.1677235007328671; "K" = 3; "M" = 14.45863566666667; "N" = 43.375907; "O" = 0.1709967203219173; "P" = 0.1709967203219173; "Q" = 17.02201682310044; "R" = 153.198151407904; "S" = 0.02868016854623097; "T" = 0.02868016854623097 }
    17 = @{ "E" = 3; "G" = 1.177290666666667; "H" = 3.531872; "I" = 0.1677235007328671; "J" = 0.1677235007328671; "K" = 3; "M" = 10.640006; "N" = 31.920018; "O" = 0.1258352566694817; "P" = 0.1258352566694817; "Q" = 12.52637975707734; "R" = 112.737417813696; "S" = 0.02110552976422433; "T" = 0.02110552976422433 }
    18 = @{ "E" = 3; "G" = 1.177290666666667; "H" = 3.531872; "I" = 0.1677235007328671; "J" = 0.1677235007328671; "K" = 3; "M" = 8.185362; "N" = 24.556086; "O" = 0.09680512663269379; "P" = 0.09680512663269379; "Q" = 9.636550285887999; "R" = 86.72895257299201; "S" = 0.01623649472772391; "T" = 0.01623649472772391 }
    19 = @{ "E" = 3; "G" = 1.177290666666667; "H" = 3.531872; "I" = 0.1677235007328671; "J" = 0.1677235007328671; "K" = 3; "M" = 3.095449666666667; "N" = 9.286349; "O" = 0.03660869207333731; "P" = 0.03660869207333731; "Q" = 3.644244001703111; "R" = 32.798196015328; "S" = 0.006140137991791698; "T" = 0.006140137991791698 }
}

foreach ($rowNum in $updates.Keys) {
    $rowData = $updates[$rowNum]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$rowNum").Value = $rowData[$col]
    }
}
